$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row (98) right after the current last row (97),
# duplicating its formatting via a row-level copy/insert so that the
# new cells pick up the exact same style records as row 97.
$ws.Rows.Item(97).Copy()
$ws.Rows.Item(98).Insert(-4121) # xlShiftDown

# Fill in the new day's data (date serial 43998 = 2020-06-16)
$ws.Cells.Item(98, 1).Value = 43998
$ws.Cells.Item(98, 2).Value = 89151
$ws.Cells.Item(98, 3).Value = 986
$ws.Cells.Item(98, 4).Value = 1503
$ws.Cells.Item(98, 5).Value = 4
$ws.Cells.Item(98, 6).Value = 7
$ws.Cells.Item(98, 7).Value = 1
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 109
$ws.Cells.Item(98, 10).Value = 0

# Extend the table range to include the new row
$table = $ws.ListObjects.Item("Tabela1")
$table.Resize($ws.Range("A1:J98"))

# Match the new active selection to the freshly appended row
$ws.Range("A98:J98").Select()
